# "Generate Report for Handback" -- mark the two localization jobs as handed
# back: in sync with en-US, fill in the Latest Target File / Latest Handback
# File / Latest Handback DateTime columns on the zh-cn and de-de sheets, and
# widen the columns that now hold the longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Same hyperlink targets already used by column A (GitHub blob links for the
# two source .md files).
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d86af8c8d19ead65f7747d5292e0be6bb130958/e2e/2be55064-692c-4b21-9ef5-06f8238408c0.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d86af8c8d19ead65f7747d5292e0be6bb130958/e2e/725e6dbe-a668-46a4-a396-f81c9bd618f1.md"
$mdName1 = "2be55064-692c-4b21-9ef5-06f8238408c0.md"
$mdName2 = "725e6dbe-a668-46a4-a396-f81c9bd618f1.md"

# Cornflower-blue (FF6495ED), same color used by the existing hyperlink style.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: the "Status" column for both locales moves from
# "In Translation" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 30
$wsOverview.Columns.Item(6).ColumnWidth = 30

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = $hyperlinkColor
$wsZh.Range("J2").Value = "2be55064-692c-4b21-9ef5-06f8238408c0.28301f2cdd04969c91342bc89a5767af14153d73.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-09 13:22:05"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = $hyperlinkColor
$wsZh.Range("J3").Value = "725e6dbe-a668-46a4-a396-f81c9bd618f1.7327dbd3aeb3ed406d21ed070587ab1fb62f47ee.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-09 13:22:05"

$wsZh.Columns.Item(3).ColumnWidth = 30
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = $hyperlinkColor
$wsDe.Range("J2").Value = "2be55064-692c-4b21-9ef5-06f8238408c0.28301f2cdd04969c91342bc89a5767af14153d73.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-09 13:22:32"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = $hyperlinkColor
$wsDe.Range("J3").Value = "725e6dbe-a668-46a4-a396-f81c9bd618f1.7327dbd3aeb3ed406d21ed070587ab1fb62f47ee.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-09 13:22:32"

$wsDe.Columns.Item(3).ColumnWidth = 30
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40
